$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text (not auto-converted numbers/dates) for Price/Volume columns, matching source data stored as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.301.04'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '1.968.75'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '248.40'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4900'
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '44.80'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2994'
$ws.Range("E9").Value = '  +3.51%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.06871'
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '19.33'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '107.41'
$ws.Range("E12").Value = '  -3.20%  '
$ws.Range("D13").Value = '1.946.09'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.07782'
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.483'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.7188'
$ws.Range("E16").Value = '  +7.11%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '288.97'
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '31.327.71'
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007786'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '5.660'
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").Value = '0.4924'
$ws.Range("E22").Value = '  +10.24%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.198.96'
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D25").Value = '1.0000'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("B26").Value = 'Chainlink'
$ws.Range("C26").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D26").Value = '6.668'
$ws.Range("E26").Value = '  +3.36%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  +6.43%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '169.54'
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '20.08'
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '2.202'
$ws.Range("E30").Value = '  +5.89%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.1071'
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '1.451'
$ws.Range("E32").Value = '  +1.22%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.881'
$ws.Range("E33").Value = '  +20.50%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.540'
$ws.Range("E34").Value = '  +9.32%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05092'
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7744'
$ws.Range("E36").Value = '  +5.25%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '1.178'
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02062'
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").Value = '2.733'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.717'
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '2.160'
$ws.Range("E41").Value = '  +6.93%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.430'
$ws.Range("E42").Value = '  +10.81%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4504'
$ws.Range("E43").Value = '  +1.21%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.8876'
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '110.13'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '73.91'
$ws.Range("E46").Value = '  +6.29%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.560'
$ws.Range("E48").Value = '  +5.03%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '996.44'
$ws.Range("E49").Value = '  +18.21%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1275'
$ws.Range("E50").Value = '  +4.00%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.419'
$ws.Range("E51").Value = '  +2.58%  '
